$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert a new parameter row ("NofInstalments" / "No Of Installments") as the
#     new row 2, push the existing SA/Rate rows down to rows 3 & 4, and rename
#     "SA" -> "Sa" on the way. NOTE: we deliberately avoid Rows.Insert() here --
#     on this host it leaves rows 2/3 aliased (writes to row 2 leak into row 3),
#     so instead every cell is written explicitly, column by column, bottom row
#     first within each column (all previously-unused row 4 cells are therefore
#     never at risk of aliasing with anything).

$ws.Range("B2").Value2 = "P0036"

$ws.Range("C4").Value2 = 3
$ws.Range("C3").Value2 = 2
$ws.Range("C2").Value2 = 1

$ws.Range("D4").Value2 = "Rate"
$ws.Range("D3").Value2 = "Sa"
$ws.Range("D2").Value2 = "NofInstalments"

$ws.Range("E4").Value2 = "float64"
$ws.Range("E3").Value2 = "float64"
$ws.Range("E2").Value2 = "int"

$ws.Range("F4").Value2 = "Null"
$ws.Range("F3").Value2 = "Null"
$ws.Range("F2").Value2 = "Null"

$ws.Range("G4").Value2 = "Null"
$ws.Range("G3").Value2 = "Null"
$ws.Range("G2").Value2 = "Null"

$ws.Range("H4").Value2 = "Null"
$ws.Range("H3").Value2 = "Null"
$ws.Range("H2").Value2 = "Null"

$ws.Range("J4").Value2 = "Rate"
$ws.Range("J3").Value2 = "Sum Assured"
$ws.Range("J2").Value2 = "No Of Installments"

$ws.Range("K4").Value2 = "y"
$ws.Range("K3").Value2 = "y"
$ws.Range("K2").Value2 = "y"

$ws.Range("Q2").Value2 = "StampDuties"

# --- Column widths: widen D (new column holding longer COLUMN_NAME values) and
#     J (Column Description, now holds "No Of Installments").
$ws.Columns.Item(4).ColumnWidth = 12.7
$ws.Columns.Item(10).ColumnWidth = 15.7

# --- Match the author's final selection/cursor position.
$ws.Range("K4").Select()
